$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Populate new data cells added in rows 85-122 ---
$ws.Range("Q85").Value = "2023年6月14日11:28:25测了一下运行速度，运行速度有点慢呀，是不是矩阵太多了。还是判断距离的那部分太多了。"
$ws.Range("R87").Value = "200秒算了一千步，说慢也不算慢，但是肯定也算不上快就是了"
$ws.Range("S87").Value = "0.2秒一步"
$ws.Range("R88").Value = "CPU占用并不高，说明是单核的。如果开个8核就可以变成几十秒一步了。"
$ws.Range("D89").Value = "一个很关键的参数就是刷新率，或者说轨道周期，总之就是过一次能用上的"
$ws.Range("R89").Value = "恐怕真的得琢磨一下什么线程池机制啊之类的"
$ws.Range("Q91").Value = "那么问题就来了，是多线程还是改矩阵呢？反正分析下来最大的时间消耗还是在node更新状态这一步。"
$ws.Range("C92").Value = "2023年6月14日14:36:37生成几何还是有问题，现在这个检测不出那种不连续的点凹进去的。"
$ws.Range("R92").Value = "看起来还是线程池比较好改，对程序的改动比较小。改矩阵会损失可读性。"
$ws.Range("D93").Value = "看来得整一个偏转角度检测，超过360度直接寄"
$ws.Range("R93").Value = "好，决定了，就改线程池。刚好相当于CPU和GPU都用起来了。"
$ws.Range("Q95").Value = "多线程改起来挺快的，但是线程是开起来了，CPU占用并没有上去。这就尴尬了。"
$ws.Range("R96").Value = "好吧，好像小有增加，但显然没有拉满。"
$ws.Range("C97").Value = "2023年6月15日14:44:02，艹，实在不行就直接C++写一个battlefiled？现在至少是掌握了可能性的。要不就是改成矩阵版本的"
$ws.Range("Q97").Value = "线程数"
$ws.Range("R97").Value = "1000步耗时"
$ws.Range("D98").Value = "要不就是进一步改进算法，node用范围定义，然后只更新周围的点"
$ws.Range("Q98").Value = 3
$ws.Range("R98").Value = 205
$ws.Range("T98").Value = "好像没什么变化"
$ws.Range("Q99").Value = 6
$ws.Range("T99").Value = "好像还是没什么变化。这相当于是计算密集型的东西了，搞来搞去都是它自己Python解释器在算，不像之前可以调用外面的东西，所以效果没什么变化？是这个道理吗。"
$ws.Range("C101").Value = "2023年6月15日16:53:52好，程序大改完事儿了。避免了大量的不必要的重复计算，代价是程序可读性下降和复杂性增加。"
$ws.Range("Q101").Value = "那意思还得攒一个C++的计算库进来？用于给你并行？"
$ws.Range("D104").Value = "好这下好了，变成6秒了"
$ws.Range("Q104").Value = "行吧，那就整一下。想必不会比JNI更加蛋疼吧"
$ws.Range("D105").Value = 6.2612535953521702
$ws.Range("E105").Value = 6.2253460884094203
$ws.Range("Q105").Value = "2023年6月14日17:21:13，实现了demo，但是改写起来恐怕挺蛋疼的。"
$ws.Range("R106").Value = "草了，有那么一瞬间甚至在怀疑，用Python写是不是错误的。"
$ws.Range("D107").Formula = "=R115/D105"
$ws.Range("Q108").Value = "那传参数这部分还挺麻烦的呀，至少传一个list进去肯定是跑不了的。"
$ws.Range("Q109").Value = "文档"
$ws.Range("R109").Value = "https://pybind11.readthedocs.io/en/stable/advanced/pycpp/object.html"
$ws.Range("Q110").Value = "不行，只能传基础的参数类型，不能传自定义的。"
$ws.Range("R111").Value = "把求距离的部分并行出去，别的仍然放在Python这头，应该是比较好的。"
$ws.Range("Q113").Value = "算是把C++的弄进去了，看一下效果"
$ws.Range("U113").Value = "CPU占比确实是增加了一些"
$ws.Range("Q114").Value = "线程数"
$ws.Range("R114").Value = "1000步耗时"
$ws.Range("Q115").Value = 3
$ws.Range("R115").Value = 190.15443372726401
$ws.Range("S115").Value = 182.33571934700001
$ws.Range("V115").Value = "算是稍微有点效果。进一步封装一下估计还能更好点儿"
$ws.Range("Q116").Value = 6
$ws.Range("R116").Value = 179.96466469764701
$ws.Range("S116").Value = 182.988606929779
$ws.Range("T116").Value = 179.21554517745901
$ws.Range("V116").Value = "和3的相比，CPU占比并没有显著的增加，说明还是寄？"
$ws.Range("Q119").Value = "好，这次是把判断也挪进去了，这样不用往node里面传UAV了，C++传回来的也直接是bool"
$ws.Range("U119").Value = "什么东西？CPU占用反而是降低了？是因为删了time，sleep？2023年6月15日14:35:57好像是，加上之后貌似恢复了"
$ws.Range("Q120").Value = "线程数"
$ws.Range("R120").Value = "1000步耗时"
$ws.Range("Q121").Value = 3
$ws.Range("R121").Formula = "=175.22634601593"
$ws.Range("S121").Value = 175.24729108810399
$ws.Range("V121").Value = "算是稍微又有一点点效果，但是显然没有实现把CPU拉满的目标。"
$ws.Range("Q122").Value = 6
$ws.Range("R122").Value = 175.37919330596901
$ws.Range("V122").Value = "总之就是也差不多，没什么前途。还得是好好优化算法。"

# --- New hyperlink cells (Q102, Q103) ---
$ws.Hyperlinks.Add($ws.Range("Q102"), "https://zhuanlan.zhihu.com/p/595502483，行吧，真有人这么干来实现多线程的。") | Out-Null
$ws.Range("Q102").Value = "https://zhuanlan.zhihu.com/p/595502483，行吧，真有人这么干来实现多线程的。"
$ws.Range("Q102").Style = "超链接"
$ws.Hyperlinks.Add($ws.Range("Q103"), "https://zhuanlan.zhihu.com/p/383572973，这个，好像靠谱。") | Out-Null
$ws.Range("Q103").Value = "https://zhuanlan.zhihu.com/p/383572973，这个，好像靠谱。"
$ws.Range("Q103").Style = "超链接"

# --- Update selection / scroll position to mirror the saved view state ---
$ws.Range("K107").Select()

# --- Page setup: A4 portrait with a linked printer-settings part ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Output "edit complete"
